# Add 2022-Q1 sheet with fund-holding data, insert it right before the
# "总计" (summary) sheet, and update the 总计 sheet with a new row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# --- Locate the summary sheet ("总计") which currently sits last ---
$summaryBefore = $wb.Worksheets.Item("总计")

# --- Create the new quarter sheet and place it right before 总计 ---
$newSheet = $wb.Worksheets.Add($summaryBefore)
$newSheet.Name = "2022-Q1"

# Re-fetch the summary sheet reference: its position shifted after the insert.
$summary = $wb.Worksheets.Item("总计")

# Copy header formatting (bold font + border + centered alignment) from an
# existing quarter sheet so the new sheet's header matches the others exactly.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Column-A style (bold font + border + centered alignment), matching the
# other quarter sheets.
$template.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# Data row. The text-like numeric codes/figures are stored as plain text
# (matching the rest of the workbook), entered with a leading apostrophe to
# force text type and then reset to the "Normal" style so no number format
# is left behind on the cell.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'162416"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "华宝港股通恒生香港35指数(LOF)"
$newSheet.Range("D2").Value = "'0.21"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'94.50"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'4.73"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.0099"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 6

# --- Update the summary sheet: insert a new row for 2022-Q1 right after the header ---
$summary.Rows.Item(2).Insert()

# Copy column-A formatting (bold font + border + centered alignment) from the
# row below (which held A2's formatting before the insert shifted it to A3).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# The inserted B2:D2 cells inherit the header's formatting by default; clear
# it so they look like ordinary data cells (matching the other data rows).
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q1"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.01

# Re-number the index column (A) for the rest of the rows, which now shifted down by one.
$lastRow = $summary.Cells.Item($summary.Rows.Count, "B").End(-4162).Row
for ($r = 3; $r -le $lastRow; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}
